$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.805.94"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "2.112.38"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.Value = "'234.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  -0.03%  "
$c = $ws.Range("D8")
$c.Value = "'57.78"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  +2.13%  "
$c = $ws.Range("D10")
$c.Value = "'0.0777"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "2.422.63"
$ws.Range("E12").Value = "  +2.54%  "
$c = $ws.Range("D13")
$c.Value = "'14.45"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$c = $ws.Range("D14")
$c.Value = "'21.22"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "
$c = $ws.Range("D15")
$c.Value = "'0.781"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "2.099.86"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "37.742.96"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("E19").Value = "  -2.34%  "
$c = $ws.Range("D20")
$c.Value = "'70.28"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("E21").Value = "  +1.83%  "
$c = $ws.Range("D22")
$c.Value = "'227.17"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("E23").Value = "  +0.07%  "
$c = $ws.Range("D24")
$c.Value = "'2.42"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +1.46%  "
$c = $ws.Range("D26")
$c.Value = "'169.33"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("E27").Value = "  +2.39%  "
$c = $ws.Range("D28")
$c.Value = "'0.134"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +6.13%  "
$ws.Range("E29").Value = "  -1.73%  "
$c = $ws.Range("D30")
$c.Value = "'19.46"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("E32").Value = "  +3.98%  "
$c = $ws.Range("D33")
$c.Value = "'2.59"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("E34").Value = "  +1.77%  "
$c = $ws.Range("D35")
$c.Value = "'4.55"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.81%  "
$c = $ws.Range("D36")
$c.Value = "'3.44"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +5.50%  "
$c = $ws.Range("D37")
$c.Value = "'1.82"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("E38").Value = "  -0.11%  "
$c = $ws.Range("D39")
$c.Value = "'5.43"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -6.03%  "
$c = $ws.Range("D40")
$c.Value = "'0.100"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +8.32%  "
$ws.Range("E41").Value = "  -0.43%  "
$c = $ws.Range("D42")
$c.Value = "'96.63"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").Value = "1.462.54"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("E45").Value = "  -0.33%  "
$c = $ws.Range("D46")
$c.Value = "'4.11"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -11.23%  "
$ws.Range("E47").Value = "  +3.35%  "
$c = $ws.Range("D48")
$c.Value = "'15.43"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "
$c = $ws.Range("D49")
$c.Value = "'3.04"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.00%  "
$c = $ws.Range("D50")
$c.Value = "'7.27"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").Value = "2.307.71"
$ws.Range("E51").Value = "  +2.57%  "
